$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two existing date headers ---
$ws.Range("A62").Value = "日期：2017.9.19 第四周 周二"
$ws.Range("A72").Value = "日期：2017.9.20 第四周 周三"

# --- Clone the previous week block (A62:D70) down to the new block (A82:D90) ---
# First copy values (and merges), then copy formats on top so styles line up with
# the template block exactly.
$ws.Range("A62:D70").Copy()
$ws.Range("A82").PasteSpecial(-4163)
$ws.Range("A62:D70").Copy()
$ws.Range("A82").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-establish the merges for the new header/footer rows (copy above only carried
# the per-cell formats, not the merge geometry).
$ws.Range("A82:D82").Merge()
$ws.Range("A90:D90").Merge()

# The new week's entries are short one-liners, so every row in the block sits at
# the sheet's single-line auto height (same as the other short header/value rows).
$ws.Rows("82:90").RowHeight = 22.5

# --- Populate the new block's content ---
$ws.Range("A82").Value = "日期：2017.9.25 第五周 周一"

$ws.Range("B84").Value = "编写数据库设计文档"
$ws.Range("B85").Value = "编写用例规约"
$ws.Range("B86").Value = "编写用例规约"
$ws.Range("B87").Value = "编写用例规约"
$ws.Range("B88").Value = "编写数据库设计文档"
$ws.Range("B89").Value = "编写用例规约"

# "完成情况" column is not filled in yet for the new week.
$ws.Range("C84").ClearContents()
$ws.Range("C85").ClearContents()
$ws.Range("C86").ClearContents()
$ws.Range("C87").ClearContents()
$ws.Range("C88").ClearContents()
$ws.Range("C89").ClearContents()
$ws.Range("D84").ClearContents()

$ws.Range("A90").Value = "总结："

# --- Update the view so the new rows are visible/selected, matching the source file ---
$ws.Application.ActiveWindow.ScrollRow = 71
$ws.Range("A90:D90").Select()
